$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-13 (A:F), after re-sorting/updating candidate statuses.
$data = @(
    @(677, "PointFive", "VP EMEA", "Jeff Ganon", "2nd Interview", 45988),
    @(677, "PointFive", "VP EMEA", "Jeremy Perlman", "2nd Interview", 45993),
    @(677, "PointFive", "VP EMEA", "Patrick Wittenberg", "3rd Interview", 45992),
    @(677, "PointFive", "VP EMEA", "Peter Reeve", "1st Interview", 45991),
    @(702, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (UK)", "Andrej Chomutovskij", "1st Interview", 45991),
    @(702, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (UK)", "Ciaran Deasy", "4th Interview", 45994),
    @(702, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (UK)", "Michel Kiffel", "4th Interview", 45993),
    @(702, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (UK)", "Simone Malekar", "1st Interview", 45991),
    @(714, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (Middle East)", "Aamer Mushtaq", "CV Sent", 45983),
    @(714, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (Middle East)", "Andrej Chomutovskij", "CV Sent", 45983),
    @(714, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (Middle East)", "JOSHUA TANNER", "1st Interview", 45991),
    @(714, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (Middle East)", "Oliver Waterman", "1st Interview", 45991)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
